# Q factor run for sg_rr_20_025 2023-12-13 17-59-26.csv data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "prominence" note in A70 to the longer, more detailed version.
$ws.Range("A70").Value = "I will choose the prominence roughly by looking at height span of roughly biggest height span noise. But it is very rough, and main check is whether code visually appears to find peaks. I may actually make the prominence a bit less that full height span of biggest height span on noise depending on how high the peaks seem to be above the noise. I will try and see what I can get away with."

# Insert two new rows before the old header row (71 and 72), pushing the
# header table (and the stray helper note below it) down by two rows.
$ws.Rows("71:72").Insert()

$ws.Range("A71").Value = "I will use same approx fsrs as above."

# Add the two new header columns for radius / radius error to the (now
# shifted) header row. (Set before A72 below so the shared-string table
# ends up in the same order as the reference workbook.)
$ws.Range("S74").Value = "radis/micrometres"
$ws.Range("T74").Value = "radius error/ micrometres"

$ws.Range("A72").Value = "I get the radius from the filename, and assume the error for all to be 0.1 micrometres."

# Fill in the new data row (75) for the sg_rr_20_025 2023-12-13 17-59-26.csv run.
$ws.Range("A75").Value = "sg_rr_20_025 2023-12-13 17-59-26.csv"
$ws.Range("B75").Value = 0.01
$ws.Range("C75").Value = 1000
$ws.Range("D75").Value = 5001
$ws.Range("E75").Value = 1530
$ws.Range("F75").Value = 1570
$ws.Range("G75").Value = 0.004
$ws.Range("H75").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I75").Value = 5
$ws.Range("J75").Value = 4.92
$ws.Range("K75").Value = 0.025354627641843101
$ws.Range("L75").Value = "yes"
$ws.Range("M75").Value = 0.15234521434913001
$ws.Range("N75").Value = 0.0160120635742305
$ws.Range("O75").Value = 11059.505619039001
$ws.Range("P75").Value = 1285.5192974044101
$ws.Range("Q75").Value = 3381794001.8318701
$ws.Range("R75").Value = 1179749480.62256
$ws.Range("S75").Value = 20
$ws.Range("T75").Value = 0.1

# Update sheet view to match the scrolled/selected state after the edit.
$ws.Application.ActiveWindow.ScrollRow = 65
$ws.Range("A76").Select()
